# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { cell address -> new value }
$updates = @{
    "展览" = @{
        "F3"  = 1339
        "F8"  = 11587
        "F9"  = 4367
        "F15" = 1089
        "F16" = 139
        "F18" = 4525
        "F22" = 11245
    }
    "全部类型" = @{
        "F3"  = 1339
        "F8"  = 11587
        "F9"  = 4367
        "F16" = 1089
        "F17" = 139
        "F19" = 4525
        "F23" = 11245
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($addr in $cellUpdates.Keys) {
        $ws.Range($addr).Value = $cellUpdates[$addr]
    }
}
